$wb = $excel.ActiveWorkbook

# Sheet "11" (the first tab) gets a new list of 4 names added in A1:A4,
# becomes the active/selected sheet (taking that role from sheet "13"),
# with the active cell left on E3.
$ws1 = $wb.Worksheets.Item("11")
$ws1.Activate()

$ws1.Range("A1").Value = "batel elbaz"
$ws1.Range("A2").Value = "shahar gavriel"
# Write A4 before A3 so the shared-string table ends up in the same
# order as the source edit (liad tzvaot = index 10, idan yontov = index 11).
$ws1.Range("A4").Value = "liad tzvaot"
$ws1.Range("A3").Value = "idan yontov"

$ws1.Range("E3").Select()
